# Stage 1 body parts: insert 6 new Key/Value rows ("bodyBacillus", etc.)
# above the existing "cellStructure*"/"motilityFlagella*" rows on the "en"
# sheet, shifting everything from row 14 down to row 20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing rows 14-21 down to 20-27, leaving rows 14-19 empty.
$ws.Rows("14:19").Insert()

# Fill the newly inserted rows. The fill order below (16, then 18-19,
# then 17) mirrors how the rows were actually populated so the shared
# string table ends up in the same order as the source edit.
$ws.Range("A14").Value = "bodyBacillus"
$ws.Range("B14").Value = "Bacillus"

$ws.Range("A15").Value = "bodyBacillusTiny"
$ws.Range("B15").Value = "Tiny Bacillus"

$ws.Range("A16").Value = "bodyCoccus"
$ws.Range("B16").Value = "Coccus"

$ws.Range("A18").Value = "bodyCoccobacillus"
$ws.Range("B18").Value = "Coccobacillus"

$ws.Range("A19").Value = "bodyCorkscrew"
$ws.Range("B19").Value = "Corkscrew"

$ws.Range("A17").Value = "bodyCoccusTiny"
$ws.Range("B17").Value = "Tiny Coccus"

# Match the new selected cell recorded in the workbook.
$ws.Range("A15").Select()
